# Add a new "Correction " column (N) to the "Card4" sheet, matching the
# pattern already used on sibling "CardN" sheets (e.g. Card1/Card2/Card7):
#   - M1 header text loses its trailing space: "Event " -> "Event"
#   - N1 gets a new header "Correction " with the same style as the other
#     header cells (bold, centered, bordered)
#   - M2:M13 (previously blank placeholder cells) become the text "nan"
#   - N2:N13 become new blank placeholder cells, like M2:M13 used to be

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card4")

# While M2:M13 are still blank, clone their (un-styled) blank format into
# the new N2:N13 cells so the new column materializes with plain formatting
# instead of picking up a stray default style.
$ws.Range("M2:M13").Copy()
$ws.Range("N2:N13").PasteSpecial(-4122)

# Clone the bold/centered/bordered header style from M1 onto the new N1
# header cell, then set the header texts.
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").Value = "Correction "
$ws.Range("M1").Value = "Event"

# Populate the existing M column's blank rows with the literal text "nan",
# matching the other sheets' placeholder convention.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}
